$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D text values that look numeric (e.g. "226.49", "1.00")
# from being auto-converted to numbers by Excel when assigned via .Value
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '39.755.35'
$ws.Range("E2").Value = '  +2.43%  '
$ws.Range("D3").Value = '2.159.55'
$ws.Range("E3").Value = '  +2.70%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '226.49'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("D7").Value = '62.82'
$ws.Range("E7").Value = '  +1.02%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").Value = '2.480.86'
$ws.Range("E13").Value = '  +2.67%  '
$ws.Range("D14").Value = '21.81'
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '0.804'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").Value = '5.51'
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("D17").Value = '2.157.53'
$ws.Range("E17").Value = '  +2.17%  '
$ws.Range("D18").Value = '39.679.13'
$ws.Range("E18").Value = '  +2.27%  '
$ws.Range("D19").Value = '71.75'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("D21").Value = '0.0₃0844'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = '229.72'
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '2.41'
$ws.Range("E24").Value = '  +2.62%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("D26").Value = '172.27'
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("D27").Value = '9.53'
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("E28").Value = '  +1.47%  '
$ws.Range("D29").Value = '1.46'
$ws.Range("E29").Value = '  +3.20%  '
$ws.Range("D30").Value = '19.78'
$ws.Range("E30").Value = '  +2.06%  '
$ws.Range("D31").Value = '2.70'
$ws.Range("E31").Value = '  +6.55%  '
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("E34").Value = '  -2.38%  '
$ws.Range("D35").Value = '6.90'
$ws.Range("E35").Value = '  -3.34%  '
$ws.Range("D36").Value = '0.0616'
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").Value = '3.77'
$ws.Range("E37").Value = '  +7.21%  '
$ws.Range("D38").Value = '2.39'
$ws.Range("E38").Value = '  +0.83%  '
$ws.Range("B39").Value = 'BinanceUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("B40").Value = 'FTXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D40").Value = '4.94'
$ws.Range("E40").Value = '  +18.92%  '
$ws.Range("D41").Value = '102.26'
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("D42").Value = '0.0228'
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("D43").Value = '17.63'
$ws.Range("E43").Value = '  -2.16%  '
$ws.Range("D44").Value = '1.510.12'
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("B47").Value = 'HuobiToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D47").Value = '2.80'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.0919'
$ws.Range("E48").Value = '  +0.72%  '
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("D50").Value = '50.21'
$ws.Range("E50").Value = '  +9.12%  '
$ws.Range("B51").Value = 'TerraClassic'
$ws.Range("C51").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D51").Value = '0.000190'
$ws.Range("E51").Value = '  +32.66%  '

# Restore default style on column D so no stray formatting/style markers remain
$dRange.Style = "Normal"
